$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (row 10) loses its special "latest" date-only
# format and reverts to the standard datetime number format used by
# the other data rows.
$ws.Cells.Item(10, 1).NumberFormat = $ws.Cells.Item(9, 1).NumberFormat

# Append the new day's results as row 11.
$ws.Cells.Item(11, 1).Value = 45751
$ws.Cells.Item(11, 2).Value = 39
$ws.Cells.Item(11, 3).Value = 37
$ws.Cells.Item(11, 4).Value = 38

# The newest row gets the special "latest" date-only format that used
# to belong to row 10.
$ws.Cells.Item(11, 1).NumberFormat = "YYYY-MM-DD"
